$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - replace existing sample values with the new configuration values
$ws.Range("A2").Value = "DCM_DIM_FACT_TEST"
$ws.Range("B2").Value = "Placement ID"
$ws.Range("C2").Value = "No"
$ws.Range("D2").Value = "VARCHAR(2000)"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "kinesso.us-east-1"
$ws.Range("G2").Value = "UM_AMEX_US"
$ws.Range("H2").Value = "GR_KINESSO"
$ws.Range("I2").Value = "UM_AMEX_US"
$ws.Range("J2").Value = "UM_AMEX_US"

# Row 3 - new row of configuration values
$ws.Range("A3").Value = "SA360_GOPRO"
$ws.Range("B3").Value = "Account, From, Campaign"
$ws.Range("C3").Value = "Yes"
$ws.Range("D3").Value = "Auto"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "kinesso.us-east-1"
$ws.Range("G3").Value = "UM_GOPRO_US"
$ws.Range("H3").Value = "GR_KINESSO"
$ws.Range("I3").Value = "UM_GOPRO_US"
$ws.Range("J3").Value = "UM_GOPRO_US"

# Update selection to match final saved state
$ws.Range("G3").Select()
